# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Commit message: "regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals"
# The only cells that actually changed in the diff are column G (header "K") values
# for rows 2-18; they are replaced with the real strikeout counts computed from
# the refreshed source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 5
    4  = 6
    5  = 3
    6  = 2
    7  = 2
    8  = 3
    9  = 1
    10 = 5
    11 = 3
    12 = 4
    13 = 1
    14 = 0
    15 = 1
    16 = 3
    17 = 5
    18 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
